$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.589.47"
$ws.Range("E2").Value = "  -2.67%  "
$ws.Range("D3").Value = "3.389.01"
$ws.Range("E3").Value = "  -2.63%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'405.82"
$ws.Range("E5").Value = "  -2.41%  "
$ws.Range("D6").Value = "'133.31"
$ws.Range("E6").Value = "  +7.70%  "
$ws.Range("D7").Value = "'0.590"
$ws.Range("E7").Value = "  -2.57%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.667"
$ws.Range("E9").Value = "  -2.89%  "
$ws.Range("E10").Value = "  -8.68%  "
$ws.Range("D11").Value = "'42.39"
$ws.Range("E11").Value = "  +2.46%  "
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "3.916.28"
$ws.Range("E13").Value = "  -2.60%  "
$ws.Range("D14").Value = "'8.38"
$ws.Range("E14").Value = "  -2.92%  "
$ws.Range("D15").Value = "'19.70"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "3.389.72"
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("D17").Value = "61.590.73"
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("D19").Value = "'10.93"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").Value = "'0.0000127"
$ws.Range("E20").Value = "  -11.27%  "
$ws.Range("D21").Value = "'3.19"
$ws.Range("E21").Value = "  -4.08%  "
$ws.Range("D22").Value = "'84.95"
$ws.Range("D23").Value = "'315.02"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "'12.75"
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").Value = "'3.13"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").Value = "'4.78"
$ws.Range("E26").Value = "  +11.02%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'29.46"
$ws.Range("E27").Value = "  -5.33%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'8.30"
$ws.Range("E28").Value = "  +5.19%  "
$ws.Range("D29").Value = "'7.57"
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").Value = "'2.66"
$ws.Range("E30").Value = "  +3.97%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "'0.170"
$ws.Range("E32").Value = "  -3.32%  "
$ws.Range("D33").Value = "'11.34"
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "'40.98"
$ws.Range("E35").Value = "  -2.14%  "
$ws.Range("D36").Value = "'0.0479"
$ws.Range("E36").Value = "  -2.33%  "
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").Value = "'3.40"
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("D40").Value = "'2.93"
$ws.Range("E40").Value = "  -3.76%  "
$ws.Range("D41").Value = "'139.64"
$ws.Range("E41").Value = "  +3.07%  "
$ws.Range("D42").Value = "'1.97"
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.124"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.295"
$ws.Range("E44").Value = "  +3.89%  "
$ws.Range("D45").Value = "'3.98"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").Value = "'16.54"
$ws.Range("E46").Value = "  -3.49%  "
$ws.Range("D47").Value = "'2.22"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").Value = "'21.34"
$ws.Range("E48").Value = "  -3.26%  "
$ws.Range("D49").Value = "2.112.84"
$ws.Range("E49").Value = "  -3.29%  "
$ws.Range("D50").Value = "'2.28"
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("D51").Value = "'1.89"
$ws.Range("E51").Value = "  -0.83%  "
